$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for row 3 (2b9c1533-...)
$wsOverview.Range("G3").Value = "2016-08-22 10:25:45"

# de-de sheet: row 3 (2b9c1533-...) Correspond Handoff Datetime (H3) shares the same
# text as Overview!G3, so it must be updated too to keep both in sync.
$wsDeDe.Range("H3").Value = "2016-08-22 10:25:45"

# zh-cn sheet: row 3 (2b9c1533-...) Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsZhCn.Range("H3").Value = "2016-08-22 10:25:41"
$wsZhCn.Range("K3").Value = "2016-08-22 10:25:59"

# de-de sheet: row 3 (2b9c1533-...) Correspond Handback DateTime (K3)
$wsDeDe.Range("K3").Value = "2016-08-22 10:26:15"
